$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.109.39"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "1.866.46"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.31"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5156"
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3763"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07158"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8898"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07599"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.871.10"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.313"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.71"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008468"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.06"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.133.60"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.029"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "2.098.76"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.49"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.459"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.95"
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.95"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.093"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.84"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.659"
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.665"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09124"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05107"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7251"
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02032"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.080"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.497"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.074"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5303"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.468"
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.73"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.280"
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1465"
$ws.Range("E45").Value = "  -3.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.908"
$ws.Range("E48").Value = "  -5.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.566"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.53"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.55"
$ws.Range("E51").Value = "  -4.90%  "
# Rows 46 and 47: PaxDollar now ranks above Decentraland in this
# snapshot, so the two rows' coin/link/price/volume content swaps,
# each with freshly-fetched Price/Volume(1h) values.
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4622"
$ws.Range("E47").Value = "  -3.76%  "
